# Update values in result_data_RandomForest worksheet (per commit "Update Name of Algo")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = -11.7917
$ws.Range("B3").Value = 6.020000000000001
$ws.Range("D3").Value = -6.798099999999994
$ws.Range("D12").Value = -7.279900000000001
$ws.Range("B14").Value = 5.845800000000001
$ws.Range("B21").Value = 9.680499999999999
$ws.Range("B23").Value = 8.877100000000004
$ws.Range("D24").Value = -7.376900000000004
$ws.Range("B25").Value = 5.008599999999998
$ws.Range("C25").Value = -13.8659
$ws.Range("D25").Value = -8.255000000000001
$ws.Range("B26").Value = 5.223300000000005
$ws.Range("C27").Value = -13.25989999999999
$ws.Range("B29").Value = 5.114400000000004
$ws.Range("C31").Value = -13.2712
$ws.Range("C39").Value = -12.73690000000001
$ws.Range("C48").Value = -11.53799999999999
$ws.Range("D50").Value = -8.163600000000001
$ws.Range("C51").Value = -11.38139999999999
$ws.Range("C52").Value = -11.3687
$ws.Range("B53").Value = 5.328600000000002
$ws.Range("D53").Value = -5.790300000000001
$ws.Range("C55").Value = -13.6892
$ws.Range("C56").Value = -12.58929999999999
$ws.Range("B57").Value = 5.118899999999997
$ws.Range("C57").Value = -13.67109999999999
$ws.Range("D57").Value = -8.142799999999996
$ws.Range("B59").Value = 4.990399999999997
$ws.Range("D61").Value = -7.946599999999998
$ws.Range("D63").Value = -7.978500000000002
$ws.Range("B69").Value = 5.256399999999994
$ws.Range("D70").Value = -7.606999999999992
$ws.Range("C73").Value = -12.8973
$ws.Range("B79").Value = 9.151800000000003
$ws.Range("B83").Value = 5.555299999999994
$ws.Range("D86").Value = -7.4891
$ws.Range("C89").Value = -10.5567
$ws.Range("C90").Value = -12.12
$ws.Range("B91").Value = 5.073299999999998
$ws.Range("C92").Value = -11.5228
$ws.Range("B93").Value = 5.829099999999998
$ws.Range("D98").Value = -8.730199999999998
$ws.Range("D100").Value = -8.328900000000004
$ws.Range("D102").Value = -7.831299999999996
